$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 6; everything currently at row 6
# and below shifts down by one (dimension grows from R103 to R104).
$ws.Rows("6:6").Insert()

# Populate the newly inserted row 6 with the new data record.
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "Femacal de La Calera"
$ws.Range("C6").Value = "Coquimbo"
$ws.Range("D6").Value = "2023-04-27"
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 100112035
$ws.Range("G6").Value = "Bruselas (repollito)"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 70
$ws.Range("K6").Value = 11500
$ws.Range("L6").Value = 12000
$ws.Range("M6").Value = 11750
$ws.Range("N6").Value = "$/malla 10 kilos"
$ws.Range("O6").Value = "Provincia de Quillota"
$ws.Range("P6").Value = 1175
$ws.Range("Q6").Value = 10
$ws.Range("R6").Value = "Hortaliza"
